$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Sender"

$senderEmail = "iacopo.depalatis@acpsystem.eu"
for ($r = 3; $r -le 13; $r++) {
    $ws.Cells.Item($r, 6).Value = $senderEmail
}
